$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "Soybean.*" measurement headers to "Mungbean.*" (row 1) ---
# Re-assigning these string values causes the now-unused "Soybean.*" shared
# strings to be dropped and the new "Mungbean.*" strings appended, which is
# exactly how the workbook's sharedStrings table was reordered in the diff.
$ws.Range("E1").Value = "Mungbean.Phenology.AccumulatedTT"
$ws.Range("F1").Value = "Mungbean.Leaf.NodeNumber"
$ws.Range("G1").Value = "Mungbean.Node.NumberError"
$ws.Range("H1").Value = "Mungbean.Leaf.BranchNumber"
$ws.Range("I1").Value = "Mungbean.Leaf.Wt"
$ws.Range("J1").Value = "Mungbean.Leaf.WtError"
$ws.Range("K1").Value = "Mungbean.Stem.Wt"
$ws.Range("L1").Value = "Mungbean.Stem.WtError"
$ws.Range("M1").Value = "Mungbean.Leaf.Area"
$ws.Range("N1").Value = "Mungbean.AboveGround.Wt"
$ws.Range("O1").Value = "Mungbean.Phenology.StartFloweringDAS"

# --- Add the new observation row (row 12) ---
$ws.Range("A12").Value = "Gatton"
$ws.Range("B12").Value = "ExtraPhenSowOctCvJade"

# D12 gets the same date format as the other Clock.Today cells (copy format
# from D2, then overwrite with the date value - 4 Jan 2022 = serial 44565).
$ws.Range("D2").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D12").Value = 44565

# C12 = D12 - D$2 (DAS since sowing), same relationship as the rows above.
$ws.Range("C12").Formula = "=D12-D`$2"
$ws.Range("C12").Style = "Normal"

$ws.Range("I12").Value = 255
$ws.Range("J12").Value = 20
$ws.Range("K12").Value = 321
$ws.Range("L12").Value = 45
$ws.Range("N12").Formula = "=K12+I12"

# --- Update the selected cell shown when the workbook is reopened ---
$ws.Range("E2").Select()
